$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phaseshift_metadata")

# Add the new data row (row 6) for the 2024-04-26 analysis run
$ws.Range("A6").Value = "2024-04-26_D_e.dat"
$ws.Range("B6").Value = "D"
$ws.Range("C6").Value = 0
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D6").Value = (Get-Date -Year 2024 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E6").Value = 2.5
$ws.Range("F6").Value = 1.8
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 202.1
$ws.Range("I6").Value = "2023-11-13_E"

# Update the active selection to match the saved view state
$ws.Range("J9").Select()
